$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.85"
$ws.Range("E2").Value = "'-5.01%"
$ws.Range("D3").Value = "'40.32"
$ws.Range("E3").Value = "'-1.78%"
$ws.Range("D4").Value = "'5.048"
$ws.Range("E4").Value = "'-3.03%"
$ws.Range("D5").Value = "'0.07414"
$ws.Range("E5").Value = "'-3.42%"
$ws.Range("D6").Value = "'1.588"
$ws.Range("E6").Value = "'-2.57%"
$ws.Range("D7").Value = "'0.9261"
$ws.Range("E7").Value = "'1.29%"
$ws.Range("E8").Value = "'-0.42%"
$ws.Range("D9").Value = "'0.1182"
$ws.Range("E9").Value = "'-3.35%"
$ws.Range("D10").Value = "'0.1750"
$ws.Range("E10").Value = "'-3.85%"
$ws.Range("D11").Value = "'0.08724"
$ws.Range("E11").Value = "'-4.25%"
$ws.Range("D12").Value = "'0.04203"
$ws.Range("E12").Value = "'-1.16%"
$ws.Range("D13").Value = "'0.1055"
$ws.Range("E13").Value = "'0.31%"
$ws.Range("D14").Value = "'0.001275"
$ws.Range("E14").Value = "'1.37%"
$ws.Range("D15").Value = "'0.005894"
$ws.Range("E15").Value = "'2.49%"
$ws.Range("D16").Value = "'3.364"
$ws.Range("E16").Value = "'0.66%"
$ws.Range("D17").Value = "'4.332"
$ws.Range("E17").Value = "'0.79%"
$ws.Range("D18").Value = "'0.3349"
$ws.Range("E18").Value = "'0.40%"
$ws.Range("D19").Value = "'7.677"
$ws.Range("E19").Value = "'3.84%"
$ws.Range("D20").Value = "'0.1363"
$ws.Range("E20").Value = "'-1.39%"
$ws.Range("D21").Value = "'0.2827"
$ws.Range("E21").Value = "'4.21%"
$ws.Range("D22").Value = "'0.03868"
$ws.Range("E22").Value = "'-3.65%"
$ws.Range("D23").Value = "'0.001298"
$ws.Range("E23").Value = "'2.78%"
$ws.Range("D24").Value = "'0.003508"
$ws.Range("E24").Value = "'-19.84%"
$ws.Range("D25").Value = "'0.0001310"
$ws.Range("D26").Value = "'0.0003757"
$ws.Range("E26").Value = "'-95.00%"
$ws.Range("D38").Value = "'0.02305"
$ws.Range("E38").Value = "'-7.82%"
$ws.Range("D39").Value = "'0.05007"
$ws.Range("E39").Value = "'-5.60%"
$ws.Range("D40").Value = "'0.007751"
$ws.Range("E40").Value = "'-1.09%"
$ws.Range("E41").Value = "'123.95%"
$ws.Range("D42").Value = "'0.1278"
$ws.Range("E42").Value = "'-2.77%"
$ws.Range("D43").Value = "'0.007427"
$ws.Range("E43").Value = "'13.54%"
$ws.Range("D44").Value = "'0.007144"
$ws.Range("E44").Value = "'-10.99%"
$ws.Range("D45").Value = "'0.3188"
$ws.Range("E45").Value = "'4.95%"
$ws.Range("D46").Value = "'0.00006722"
$ws.Range("E46").Value = "'0.11%"
$ws.Range("D47").Value = "'0.00000000756"
$ws.Range("E47").Value = "'0.69%"
$ws.Range("E48").Value = "'-28.08%"
$ws.Range("D49").Value = "'0.004239"
$ws.Range("E49").Value = "'36.62%"
$ws.Range("D50").Value = "'0.00002116"
$ws.Range("E50").Value = "'0.69%"
$ws.Range("D51").Value = "'0.0002016"
$ws.Range("E51").Value = "'0.69%"
